$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 <= former Row 8 content; Row 8 <= former Row 7 content
$ws.Range("A7").Value = 131198466
$ws.Range("B7").Value = 79834
$ws.Range("E7").Value = 229821
$ws.Range("F7").Value = "Vedflamlav"
$ws.Range("G7").Value = "Ramboldia elabens"
$ws.Range("H7").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 466092
$ws.Range("R7").Value = 6789074
$ws.Range("S7").Value = 10

$ws.Range("A8").Value = 131197802
$ws.Range("B8").Value = 57881
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "äldre spår"
$ws.Range("Q8").Value = 465938
$ws.Range("R8").Value = 6789021
$ws.Range("S8").Value = 10

# Row 17 <= former Row 18 content; Row 18 <= former Row 17 content
$ws.Range("A17").Value = 131198252
$ws.Range("B17").Value = 79244
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("M17").ClearContents()
$ws.Range("Q17").Value = 466111
$ws.Range("R17").Value = 6789063
$ws.Range("S17").Value = 10

$ws.Range("A18").Value = 131198195
$ws.Range("B18").Value = 57881
$ws.Range("E18").Value = 100049
$ws.Range("F18").Value = "Spillkråka"
$ws.Range("G18").Value = "Dryocopus martius"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("M18").Value = "färska spår"
$ws.Range("Q18").Value = 466050
$ws.Range("R18").Value = 6788971
$ws.Range("S18").Value = 10

# Row 28 <= former Row 29 content; Row 29 <= former Row 28 content
$ws.Range("A28").Value = 131199091
$ws.Range("B28").Value = 57884
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = "Tretåig hackspett"
$ws.Range("G28").Value = "Picoides tridactylus"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("M28").Value = "färska spår"
$ws.Range("Q28").Value = 466114
$ws.Range("R28").Value = 6788962
$ws.Range("S28").Value = 10
$ws.Range("AC28").Value = "2 bild, tall med gran till vänster."

$ws.Range("A29").Value = 131198231
$ws.Range("B29").Value = 79244
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("M29").ClearContents()
$ws.Range("Q29").Value = 466052
$ws.Range("R29").Value = 6789006
$ws.Range("S29").Value = 50
$ws.Range("AC29").Value = "Rikligt till måttligt i en radie av ca 50 meter  1 bild gren i förgrund"
